# Updates crypto price/volume data, applying the diff between the
# previous scrape and the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.583.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.57%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.786.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.69%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.87%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.555'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.66%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.84'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.74%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.280'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.88%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0676'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.08%  '

# Row 11
$ws.Range("E11").Value = '  +1.42%  '

# Row 12
$ws.Range("E12").Value = '  +0.78%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +12.13%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.784.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.62%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.591.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.63%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.630'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.71%  '

# Row 17
$ws.Range("E17").Value = '  +2.72%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.10%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0771'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.42%  '

# Row 21
$ws.Range("E21").Value = '  -0.10%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.38%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.94%  '

# Row 24
$ws.Range("E24").Value = '  +0.32%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.42%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.98%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.29%  '

# Row 28
$ws.Range("E28").Value = '  -0.18%  '

# Row 29
$ws.Range("E29").Value = '  -0.03%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.66%  '

# Row 31
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$ws.Range("E32").Value = '  -0.22%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.56'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.37%  '

# Row 34
$ws.Range("E34").Value = '  +0.92%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.440.87'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.78%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.10%  '

# Row 37
$ws.Range("E37").Value = '  +2.18%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.629'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.00%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '82.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.27%  '

# Row 40
$ws.Range("E40").Value = '  +4.51%  '

# Row 41
$ws.Range("E41").Value = '  +0.65%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.900'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.73%  '

# Row 43
$ws.Range("E43").Value = '  -1.00%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0505'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.23%  '

# Row 45
$ws.Range("E45").Value = '  +2.60%  '

# Row 46
$ws.Range("E46").Value = '  -2.21%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.943.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.74%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.64'
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.01%  '

# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.67%  '

# Row 51
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0121'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.06%  '
